$d = $word.ActiveDocument
$paras = $d.Paragraphs

# 1) Dato: cell -> append date
$p20 = $paras.Item(20)
$p20.Range.InsertAfter("17/4-16")

# 2) Deltagere: cell -> append name
$p21 = $d.Paragraphs.Item(21)
$p21.Range.InsertAfter("Tonni")

# 3) Begivenhed/-er cell -> insert sentence before bookmark, then add empty paragraph after
$p24 = $d.Paragraphs.Item(24)
$p24.Range.InsertBefore("Rettet udkast til idb samt bdd for styreboks til at reflektere ændringen i hvordan vi anvender vores LCD skræm")

$p24b = $d.Paragraphs.Item(24)
$endRange = $p24b.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "all done"
